$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing JSON Command values (text already present in shared strings;
# underlying shared-string index shifts automatically when the table is rebuilt on save)
$ws.Range("D4").Value = "console-diagnostics"
$ws.Range("D5").Value = "console-recovery-clear-crc-error-count"

# Remove the set-date-time / get-date-time command mappings
$ws.Range("D35").Value = "N/A"

# Add/replace the JSON Command column values for the remaining rows, in the same
# order the strings were originally authored so the rebuilt shared string table
# matches the canonical ordering.
$ws.Range("D26").Value = "clear-alarm-thresholds"
$ws.Range("D25").Value = "clear-console-archive"
$ws.Range("D27").Value = "clear-calibration-offsets"
$ws.Range("D28").Value = "clear-graph-points"
$ws.Range("D32").Value = "clear-active-alarms"
$ws.Range("D33").Value = "clear-current-data"
$ws.Range("D37").Value = "query-console-time"
$ws.Range("D38").Value = "update-archive-period"
$ws.Range("D39").Value = "stop-archiving"
$ws.Range("D40").Value = "start-archiving"
$ws.Range("D12").Value = "put-year-rain"
$ws.Range("D13").Value = "put-year-et"
$ws.Range("D29").Value = "clear-cumulative-values"
$ws.Range("D30").Value = "clear-high-values, clear-highs"
$ws.Range("D31").Value = "clear-low-values, clear-lows"
$ws.Range("D11").Value = "query-hilows"

# Update the selection to match the saved view state
[void]$ws.Range("D7").Select()
